# Weekly fruit/vegetable price update: the Fecha/Volumen/Precio columns
# (D, J, K, L, M, P) for each data row are reassigned among the existing
# rows (rows 2-46) to reflect the new weekly logic. Columns A, B, C, E,
# F, G, H, I, N, O, Q, R stay the same for every row since they are
# constant market/product metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of destination row -> source row (values currently found at the
# source row get copied into the destination row).
$rowMap = @{
    2 = 32;
    3 = 5;
    4 = 14;
    5 = 34;
    6 = 43;
    7 = 37;
    8 = 11;
    9 = 18;
    10 = 30;
    11 = 44;
    12 = 23;
    13 = 19;
    14 = 27;
    15 = 26;
    16 = 46;
    17 = 45;
    18 = 17;
    19 = 12;
    20 = 2;
    21 = 35;
    22 = 38;
    23 = 28;
    24 = 22;
    25 = 8;
    26 = 21;
    27 = 7;
    28 = 36;
    29 = 25;
    30 = 31;
    31 = 4;
    32 = 33;
    33 = 39;
    34 = 40;
    35 = 3;
    36 = 6;
    37 = 9;
    38 = 16;
    39 = 15;
    40 = 13;
    41 = 41;
    42 = 29;
    43 = 42;
    44 = 20;
    45 = 24;
    46 = 10
}

# Columns touched by the reassignment (1-based): D=4, J=10, K=11, L=12, M=13, P=16
$cols = @(4, 10, 11, 12, 13, 16)

# First snapshot the current values of every touched cell (rows 2-46),
# using Value2 so dates come back as raw serial numbers rather than
# being coerced/reformatted.
$snapshot = @{}
for ($r = 2; $r -le 46; $r++) {
    foreach ($c in $cols) {
        $key = "$r-$c"
        $snapshot[$key] = $ws.Cells.Item($r, $c).Value2
    }
}

# Now write the snapshot values back out according to the row map.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $srcKey = "$srcRow-$c"
        $ws.Cells.Item($destRow, $c).Value2 = $snapshot[$srcKey]
    }
}
